$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "CheckSelf"
$ws.Range("B4").Value = "Dan Brown"
$ws.Range("C4").Value = "Information about example account is shown"
$ws.Range("D4").Value = "Error Message"
$ws.Range("E4").Value = 'Error message says "Too few arguments"'

$ws.Range("C4:E4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 28.8

$ws.Range("E4").Select()
